$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("627:627").Insert()

$ws.Range("A627").Value = 5
$ws.Range("B627").Value = "Macroferia Regional de Talca"
$ws.Range("C627").Value = "Maule"
$ws.Range("D627").Value = 45223
$ws.Range("E627").Value = 7
$ws.Range("F627").Value = 100114014
$ws.Range("G627").Value = "Betarraga"
$ws.Range("H627").Value = "Sin especificar"
$ws.Range("I627").Value = "Primera"
$ws.Range("J627").Value = 4000
$ws.Range("K627").Value = 500
$ws.Range("L627").Value = 500
$ws.Range("M627").Value = 500
$ws.Range("N627").Value = "`$/paquete 5 unidades"
$ws.Range("O627").Value = "Región del Maule"
$ws.Range("P627").Value = 100
$ws.Range("Q627").Value = 5
$ws.Range("R627").Value = "Hortaliza"
